# The sheet "Hortaliza, Feria Lagunitas de Puerto Montt - Papa" tracks weekly
# price observations. This commit adds one new weekly observation, inserted
# as a new row 777 (most recent date), pushing the existing rows 777-794
# down to 778-795.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 777, shifting rows 777:794 down to 778:795.
$ws.Rows.Item(777).Insert()

# Populate the new row 777 with the new weekly observation.
$ws.Cells.Item(777, 1).Value  = 4
$ws.Cells.Item(777, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(777, 3).Value  = "Los Lagos"
$ws.Cells.Item(777, 4).Value  = 45239
$ws.Cells.Item(777, 5).Value  = 10
$ws.Cells.Item(777, 6).Value  = 100114001
$ws.Cells.Item(777, 7).Value  = "Papa"
$ws.Cells.Item(777, 8).Value  = "Patagonia"
$ws.Cells.Item(777, 9).Value  = "1a (guarda)"
$ws.Cells.Item(777, 10).Value = 300
$ws.Cells.Item(777, 11).Value = 29000
$ws.Cells.Item(777, 12).Value = 30000
$ws.Cells.Item(777, 13).Value = 29500
$ws.Cells.Item(777, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(777, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(777, 16).Value = 1180
$ws.Cells.Item(777, 17).Value = 25
$ws.Cells.Item(777, 18).Value = "Hortaliza"

# Ensure the style (date number format) on D777 matches the other date cells
# in this column (the Insert() above should already have copied it from the
# row above, but set it explicitly in case the row below's style differs).
$ws.Cells.Item(777, 4).NumberFormat = $ws.Cells.Item(778, 4).NumberFormat
